$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'59.001.11"
$ws.Range('E2').Value = '  +2.65%  '
$ws.Range('D3').Value = "'2.519.74"
$ws.Range('E3').Value = '  +3.72%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = "'533.68"
$ws.Range('E5').Value = '  +5.99%  '
$ws.Range('D6').Value = "'134.03"
$ws.Range('E6').Value = '  +4.51%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'0.567"
$ws.Range('E8').Value = '  +3.28%  '
$ws.Range('D9').Value = "'2.517.66"
$ws.Range('E9').Value = '  +3.13%  '
$ws.Range('D10').Value = "'0.0995"
$ws.Range('E10').Value = '  +4.95%  '
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').Value = "'5.25"
$ws.Range('E12').Value = '  +1.74%  '
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').Value = "'2.957.29"
$ws.Range('E14').Value = '  +3.30%  '
$ws.Range('D15').Value = "'58.900.47"
$ws.Range('E15').Value = '  +2.54%  '
$ws.Range('E16').Value = '  +3.19%  '
$ws.Range('E17').Value = '  +3.56%  '
$ws.Range('D18').Value = "'2.504.01"
$ws.Range('E18').Value = '  +2.75%  '
$ws.Range('D19').Value = "'10.65"
$ws.Range('E19').Value = '  +2.57%  '
$ws.Range('E20').Value = '  +4.18%  '
$ws.Range('D21').Value = "'321.17"
$ws.Range('E21').Value = '  +2.54%  '
$ws.Range('D22').Value = "'6.23"
$ws.Range('E22').Value = '  +10.04%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = "'65.79"
$ws.Range('E24').Value = '  +4.19%  '
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('D26').Value = "'0.997"
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('E28').Value = '  +4.70%  '
$ws.Range('E29').Value = '  +6.21%  '
$ws.Range('D30').Value = "'173.32"
$ws.Range('E30').Value = '  +2.05%  '
$ws.Range('E31').Value = '  +5.52%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').Value = "'6.36"
$ws.Range('E32').Value = '  +3.29%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Value = "'1.20"
$ws.Range('E33').Value = '  +7.22%  '
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('D35').Value = "'0.999"
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('D36').Value = "'18.17"
$ws.Range('E36').Value = '  +2.87%  '
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('D38').Value = "'3.95"
$ws.Range('E38').Value = '  +1.56%  '
$ws.Range('E39').Value = '  +5.22%  '
$ws.Range('B40').Value = 'SuiNetwork'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D40').Value = "'0.824"
$ws.Range('E40').Value = '  +10.34%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = "'36.70"
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = "'3.49"
$ws.Range('E42').Value = '  +4.32%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = "'277.32"
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = "'5.09"
$ws.Range('E44').Value = '  +5.39%  '
$ws.Range('D45').Value = "'131.42"
$ws.Range('E45').Value = '  +10.68%  '
$ws.Range('D46').Value = "'0.592"
$ws.Range('E46').Value = '  +2.60%  '
$ws.Range('D47').Value = "'0.0936"
$ws.Range('E47').Value = '  +3.01%  '
$ws.Range('D48').Value = "'0.0511"
$ws.Range('E48').Value = '  +5.88%  '
$ws.Range('D49').Value = "'0.0220"
$ws.Range('E49').Value = '  +5.57%  '
$ws.Range('D50').Value = "'17.08"
$ws.Range('E50').Value = '  +3.60%  '
$ws.Range('D51').Value = "'1.752.96"
$ws.Range('E51').Value = '  +3.43%  '

Write-Host "Applied 100 cell updates"
